# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the files are
# now "Ready for handoff" instead of "In Translation", and refreshes the
# handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column: "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime refresh
$overview.Range("G2").Value = "2016-08-15 22:36:42"
$dede.Range("H2").Value     = "2016-08-15 22:36:42"
$zhcn.Range("H2").Value     = "2016-08-15 22:36:36"

# Column widths widen (auto-fit-like) to accommodate the longer
# "Ready for handoff" text that now lives in these status columns.
$overview.Range("E:F").ColumnWidth = 16.4
$zhcn.Range("C:C").ColumnWidth     = 16.4
$dede.Range("C:C").ColumnWidth     = 16.4
